# "avg comment with array" - add the LeetCode 643 (Maximum Average Subarray I)
# entry as row 13 of the "数组" (Array) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("数组")

# The three new strings contain multi-line Chinese text; they are carried
# here as base64-encoded UTF-8 to dodge any PowerShell quoting issues and
# decoded back to plain text before being written into the cells.
$problemB64 = "57uZ5a6aIG4g5Liq5pW05pWw77yM5om+5Ye65bmz5Z2H5pWw5pyA5aSn5LiU6ZW/5bqm5Li6IGsg55qE6L+e57ut5a2Q5pWw57uE77yM5bm26L6T5Ye66K+l5pyA5aSn5bmz5Z2H5pWw44CCIAog56S65L6LIDE6CiDovpPlhaU6IFsxLDEyLC01LC02LDUwLDNdLCBrID0gNArovpPlh7o6IDEyLjc1Cuino+mHijog5pyA5aSn5bmz5Z2H5pWwICgxMi01LTYrNTApLzQgPSA1MS80ID0gMTIuNzUKIOazqOaEjzoKIDEgPD0gayA8PSBuIDw9IDMwLDAwMOOAggog5omA57uZ5pWw5o2u6IyD5Zu0IFstMTAsMDAw77yMMTAsMDAwXeOAgiAKIFJlbGF0ZWQgVG9waWNzIOaVsOe7hA=="
$methodB64  = "MSDorqHnrpfmlbDnu4TkuK1r5Liq6L+e57ut5pWw5a2X55qE57Sv5Yqg5ZKM5pyA5aSnCjIg5ZCR5YmN56e75Yqo5LiA5Liq5L2N572u77yM5LiN5Y+v5Lul5bCGa+S4quaVsOWtl+e0r+WKoO+8jOatpOaWueahiOaViOeOh+W+iOS9ju+8jOWPr+S7peWwhuW9k+WJjeeahOWSjOWOu+WktOWKoOWwvu+8jOiOt+WPlm1heAozIOazqOaEj++8muS9v+eUqOWJjWvkuKrmlbDlrZfliJ3lp4vljJbvvIzov63ku6PojIPlm7TmmK9bMSxsZW4gLSBrICsgMSk="
$keywordsB64 = "56qX5Y+j566X5rOVCue0r+WKoArorqHnrpfmnIDlpKflgLw="

$enc = [System.Text.Encoding]::UTF8
$problemText  = $enc.GetString([System.Convert]::FromBase64String($problemB64))
$methodText   = $enc.GetString([System.Convert]::FromBase64String($methodB64))
$keywordsText = $enc.GetString([System.Convert]::FromBase64String($keywordsB64))

# New row 13: No.=12, leetcode=643, 题目/解题方法/解题关键词 = the three strings
# above, 时间复杂度=O(N), 空间复杂度=O(1) (same shared strings already used by
# other rows on this sheet).
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 643
$ws.Cells.Item(13, 3).Value = $problemText
$ws.Cells.Item(13, 4).Value = $methodText
$ws.Cells.Item(13, 5).Value = $keywordsText
$ws.Cells.Item(13, 6).Value = "O(N)"
$ws.Cells.Item(13, 7).Value = "O(1)"

# Match the row's wrapped-text height used by the rest of the sheet.
$ws.Rows.Item(13).RowHeight = 220

# Move the active selection from D16 to D14.
[void]$ws.Range("D14").Select()
